$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 571.25
$ws.Range("I11").Value = 571.25
$ws.Range("K11").Value = 571.25
$ws.Range("M11").Value = -431.25

$ws.Range("H18").Value = 12462.6
$ws.Range("I18").Value = 14443
$ws.Range("K18").Value = 14443
$ws.Range("M18").Value = -14159

$ws.Range("H62").Value = 5154.4165
$ws.Range("I62").Value = 4009
$ws.Range("J62").Value = 8590.666999999999
$ws.Range("K62").Value = 4009
$ws.Range("L62").Value = 8590.666999999999
$ws.Range("M62").Value = -3385
$ws.Range("N62").Value = -9838.666999999999

$ws.Range("H65").Value = 5154.4165
$ws.Range("I65").Value = 4009
$ws.Range("J65").Value = 8590.666999999999
$ws.Range("K65").Value = 20045
$ws.Range("L65").Value = 42953.335
$ws.Range("M65").Value = -16925
$ws.Range("N65").Value = -49193.335

$ws.Range("H116").Value = 3459.7334
$ws.Range("J116").Value = 4391.3335
$ws.Range("L116").Value = 4391.3335
$ws.Range("N116").Value = -11275.3335

$ws.Range("H138").Value = 2316.434
$ws.Range("I138").Value = 1382.7222
$ws.Range("K138").Value = 4148.1666
$ws.Range("M138").Value = 991.8334000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 998.0851
$ws.Range("I2").Value = 936.5714
$ws.Range("J2").Value = 1514.8
$ws.Range("K2").Value = 936.5714
$ws.Range("L2").Value = 1514.8
$ws.Range("M2").Value = -823.5714
$ws.Range("N2").Value = -1740.8

$ws.Range("H32").Value = 6680.263
$ws.Range("I32").Value = 6054.4707
$ws.Range("K32").Value = 6054.4707
$ws.Range("M32").Value = -5767.4707

$ws.Range("H45").Value = 8695.474
$ws.Range("I45").Value = 15001.75
$ws.Range("J45").Value = 4109.091
$ws.Range("K45").Value = 15001.75
$ws.Range("L45").Value = 4109.091
$ws.Range("M45").Value = -14624.75
$ws.Range("N45").Value = -4863.091

$ws.Range("H69").Value = 192459
$ws.Range("J69").Value = 192459
$ws.Range("L69").Value = 192459
$ws.Range("N69").Value = -193957

$ws.Range("H72").Value = 192459
$ws.Range("J72").Value = 192459
$ws.Range("L72").Value = 577377
$ws.Range("N72").Value = -584865

$ws.Range("H97").Value = 3250
$ws.Range("I97").Value = 2900
$ws.Range("K97").Value = 2900
$ws.Range("M97").Value = -2404

$ws.Range("H110").Value = 2535.625
$ws.Range("I110").Value = 1542.8182
$ws.Range("J110").Value = 4719.8
$ws.Range("K110").Value = 1542.8182
$ws.Range("L110").Value = 4719.8
$ws.Range("M110").Value = 502.1818000000001
$ws.Range("N110").Value = -8809.799999999999

$ws.Range("H116").Value = 998.0851
$ws.Range("I116").Value = 936.5714
$ws.Range("J116").Value = 1514.8
$ws.Range("K116").Value = 936.5714
$ws.Range("L116").Value = 1514.8
$ws.Range("M116").Value = 1357.4286
$ws.Range("N116").Value = -6102.8

$ws.Range("H122").Value = 1883.7333
$ws.Range("I122").Value = 1380
$ws.Range("K122").Value = 4140
$ws.Range("M122").Value = -1690

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 998.0851
$ws.Range("I3").Value = 936.5714
$ws.Range("J3").Value = 1514.8
$ws.Range("K3").Value = 936.5714
$ws.Range("L3").Value = 1514.8
$ws.Range("M3").Value = -822.5714
$ws.Range("N3").Value = -1742.8

$ws.Range("H86").Value = 2802.3684
$ws.Range("I86").Value = 2802.3684
$ws.Range("K86").Value = 2802.3684
$ws.Range("M86").Value = -1679.3684

$ws.Range("H89").Value = 2802.3684
$ws.Range("I89").Value = 2802.3684
$ws.Range("K89").Value = 14011.842
$ws.Range("M89").Value = -8395.841999999999

$ws.Range("H99").Value = 7928.968
$ws.Range("I99").Value = 8375.793
$ws.Range("K99").Value = 8375.793
$ws.Range("M99").Value = -6877.793

$ws.Range("H107").Value = 95458
$ws.Range("I107").Value = 250382.5
$ws.Range("J107").Value = 6929.7144
$ws.Range("K107").Value = 250382.5
$ws.Range("L107").Value = 6929.7144
$ws.Range("M107").Value = -248462.5
$ws.Range("N107").Value = -10769.7144

$ws.Range("H134").Value = 939.95
$ws.Range("I134").Value = 701.2973
$ws.Range("J134").Value = 3883.3333
$ws.Range("K134").Value = 2103.8919
$ws.Range("L134").Value = 11649.9999
$ws.Range("M134").Value = 431.1081000000004
$ws.Range("N134").Value = -16719.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 2846.7693
$ws.Range("I25").Value = 3000.6667
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 3000.6667
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -2826.6667
$ws.Range("N25").Value = -1348

$ws.Range("H116").Value = 600000
$ws.Range("J116").Value = 600000
$ws.Range("L116").Value = 600000
$ws.Range("N116").Value = -609178

$ws.Range("H122").Value = 146129.86
$ws.Range("I122").Value = 169931.5
$ws.Range("K122").Value = 509794.5
$ws.Range("M122").Value = -507344.5

$ws.Range("H132").Value = 1933.7858
$ws.Range("I132").Value = 1929.2084
$ws.Range("K132").Value = 5787.6252
$ws.Range("M132").Value = -3257.6252

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2565.9524
$ws.Range("J131").Value = 2713.0667
$ws.Range("L131").Value = 8139.2001
$ws.Range("N131").Value = -18219.2001

$ws.Range("H132").Value = 1188.7778
$ws.Range("I132").Value = 1449.75
$ws.Range("J132").Value = 980
$ws.Range("K132").Value = 13047.75
$ws.Range("L132").Value = 8820
$ws.Range("M132").Value = -10517.75
$ws.Range("N132").Value = -13880

$ws.Range("H137").Value = 2780455.2
$ws.Range("J137").Value = 3290.5417
$ws.Range("L137").Value = 9871.625100000001
$ws.Range("N137").Value = -20071.6251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 5200
$ws.Range("I29").Value = 6900
$ws.Range("J29").Value = 3500
$ws.Range("K29").Value = 6900
$ws.Range("L29").Value = 3500
$ws.Range("M29").Value = -6610
$ws.Range("N29").Value = -4080

$ws.Range("H70").Value = 5327
$ws.Range("I70").Value = 5299.6665
$ws.Range("K70").Value = 5299.6665
$ws.Range("M70").Value = -5029.6665

$ws.Range("H73").Value = 5327
$ws.Range("I73").Value = 5299.6665
$ws.Range("K73").Value = 5299.6665
$ws.Range("M73").Value = -4363.6665

$ws.Range("H122").Value = 1981.742
$ws.Range("I122").Value = 1721.6364
$ws.Range("J122").Value = 2617.5557
$ws.Range("K122").Value = 5164.9092
$ws.Range("L122").Value = 7852.6671
$ws.Range("M122").Value = -2714.9092
$ws.Range("N122").Value = -12752.6671

$ws.Range("H126").Value = 2734.2856
$ws.Range("I126").Value = 2889.2
$ws.Range("K126").Value = 8667.599999999999
$ws.Range("M126").Value = -6197.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 7788.3335
$ws.Range("I58").Value = 6682.5
$ws.Range("K58").Value = 6682.5
$ws.Range("M58").Value = -6422.5

$ws.Range("H70").Value = 45000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 45000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 45000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -45540

$ws.Range("H73").Value = 45000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 45000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 45000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -46872

$ws.Range("H108").Value = 395207.66
$ws.Range("J108").Value = 395207.66
$ws.Range("L108").Value = 395207.66
$ws.Range("N108").Value = -402887.66

$ws.Range("H132").Value = 4399.946
$ws.Range("I132").Value = 3815.7036
$ws.Range("K132").Value = 11447.1108
$ws.Range("M132").Value = -8917.110799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 13000
$ws.Range("J21").Value = 13000
$ws.Range("L21").Value = 13000
$ws.Range("N21").Value = -13470

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H35").Value = 13000
$ws.Range("J35").Value = 13000
$ws.Range("L35").Value = 13000
$ws.Range("N35").Value = -13580

$ws.Range("H113").Value = 272
$ws.Range("J113").Value = 303.7143
$ws.Range("L113").Value = 911.1428999999999
$ws.Range("N113").Value = -5251.1429

$ws.Range("H126").Value = 1818.7727
$ws.Range("I126").Value = 1795.4736
$ws.Range("K126").Value = 5386.4208
$ws.Range("M126").Value = -2916.4208
